$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "TestComplete14_"
$ws.Range("D4").Value = "TEST02-PC"
$ws.Range("D5").Value = "TestExecute-PC"
$ws.Range("D6").Value = "TestComplete14_"
